# Re-process the metadata with the newly curated dimensions:
# the "provincia", "hectareas-en-tierras-labradas-de-secano-otros-cultivos"
# and "hectareas-en-tierras-labradas-con-otros-cultivos" columns move from
# being sdmx/iaest "dimension" columns to "measure" columns, so their
# metadata rows (concept URI, dim/medida marker, datatype) are updated and
# their now-unused dimension-mapping file references (row 5) are cleared.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column F: provincia
$ws.Range("F2").Value = "iaest-measure:provincia"
$ws.Range("F3").Value = "medida"
$ws.Range("F4").Value = "xsd:int"

# Column R: hectareas-en-tierras-labradas-de-secano-otros-cultivos
$ws.Range("R2").Value = "iaest-measure:hectareas-en-tierras-labradas-de-secano-otros-cultivos"
$ws.Range("R3").Value = "medida"
$ws.Range("R4").Value = "xsd:int"
$ws.Range("R5").ClearContents()

# Column Y: hectareas-en-tierras-labradas-con-otros-cultivos
$ws.Range("Y2").Value = "iaest-measure:hectareas-en-tierras-labradas-con-otros-cultivos"
$ws.Range("Y3").Value = "medida"
$ws.Range("Y4").Value = "xsd:int"
$ws.Range("Y5").ClearContents()
